$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "2 (PD2)"
$ws.Range("D3").Value = "0 (PD0)"
$ws.Range("D4").Value = "1 (PD1)"
$ws.Range("D5").Value = "3 (PD3)"
$ws.Range("D6").Value = "4 (PD4)"
$ws.Range("D7").Value = "17 (PC3)"
$ws.Range("D8").Value = "16 (PC2)"

$ws.Range("D9").Select()
